# Auto-generated edit script: updates cached market price values on the
# ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

$ws.Range("H133").Value = 50339.2
$ws.Range("J133").Value = 50339.2
$ws.Range("L133").Value = 50339.2
$ws.Range("N133").Value = -60459.2

$ws.Range("H134").Value = 59573.332
$ws.Range("J134").Value = 59573.332
$ws.Range("L134").Value = 59573.332
$ws.Range("N134").Value = -69713.33199999999

$ws.Range("H137").Value = 1811.3422
$ws.Range("I137").Value = 1460.9333
$ws.Range("J137").Value = 3125.375
$ws.Range("K137").Value = 4382.7999
$ws.Range("L137").Value = 9376.125
$ws.Range("M137").Value = -1832.7999
$ws.Range("N137").Value = -14476.125

$ws.Range("H139").Value = 48108.57
$ws.Range("J139").Value = 48108.57
$ws.Range("L139").Value = 48108.57
$ws.Range("N139").Value = -58388.57

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10639791
$ws.Range("I74").Value = 1265.7142
$ws.Range("J74").Value = 26317618
$ws.Range("K74").Value = 1265.7142
$ws.Range("L74").Value = 26317618
$ws.Range("M74").Value = -391.7141999999999
$ws.Range("N74").Value = -26319366

$ws.Range("H77").Value = 10639791
$ws.Range("I77").Value = 1265.7142
$ws.Range("J77").Value = 26317618
$ws.Range("K77").Value = 6328.571
$ws.Range("L77").Value = 131588090
$ws.Range("M77").Value = -1960.571
$ws.Range("N77").Value = -131596826

$ws.Range("H133").Value = 28475
$ws.Range("J133").Value = 28475
$ws.Range("L133").Value = 28475
$ws.Range("N133").Value = -33535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 51697
$ws.Range("J132").Value = 51697
$ws.Range("L132").Value = 51697
$ws.Range("N132").Value = -61817

$ws.Range("H134").Value = 26476.61
$ws.Range("I134").Value = 5300.4863
$ws.Range("J134").Value = 113534
$ws.Range("K134").Value = 15901.4589
$ws.Range("L134").Value = 340602
$ws.Range("M134").Value = -13366.4589
$ws.Range("N134").Value = -345672

$ws.Range("H135").Value = 75545
$ws.Range("J135").Value = 75545
$ws.Range("L135").Value = 75545
$ws.Range("N135").Value = -85685

$ws.Range("H138").Value = 47255.555
$ws.Range("J138").Value = 47255.555
$ws.Range("L138").Value = 47255.555
$ws.Range("N138").Value = -57535.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10523461
$ws.Range("I31").Value = 1933
$ws.Range("J31").Value = 19096558
$ws.Range("K31").Value = 1933
$ws.Range("L31").Value = 19096558
$ws.Range("M31").Value = -1638
$ws.Range("N31").Value = -19097148

$ws.Range("H34").Value = 10523461
$ws.Range("I34").Value = 1933
$ws.Range("J34").Value = 19096558
$ws.Range("K34").Value = 1933
$ws.Range("L34").Value = 19096558
$ws.Range("M34").Value = -1731
$ws.Range("N34").Value = -19096962

$ws.Range("H58").Value = 4650806
$ws.Range("I58").Value = 7247816
$ws.Range("J58").Value = 668723.9
$ws.Range("K58").Value = 7247816
$ws.Range("L58").Value = 668723.9
$ws.Range("M58").Value = -7247613
$ws.Range("N58").Value = -669129.9

$ws.Range("H94").Value = 2446.6086
$ws.Range("J94").Value = 2490.9412
$ws.Range("L94").Value = 2490.9412
$ws.Range("N94").Value = -3392.9412

$ws.Range("H135").Value = 38530.77
$ws.Range("J135").Value = 39658.332
$ws.Range("L135").Value = 39658.332
$ws.Range("N135").Value = -49798.332

$ws.Range("H136").Value = 4650806
$ws.Range("I136").Value = 7247816
$ws.Range("J136").Value = 668723.9
$ws.Range("K136").Value = 21743448
$ws.Range("L136").Value = 2006171.7
$ws.Range("M136").Value = -21740898
$ws.Range("N136").Value = -2011271.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2922.5557
$ws.Range("I102").Value = 2793.2856
$ws.Range("J102").Value = 3375
$ws.Range("K102").Value = 2793.2856
$ws.Range("L102").Value = 3375
$ws.Range("M102").Value = -1171.2856
$ws.Range("N102").Value = -6619

$ws.Range("H126").Value = 16453.846
$ws.Range("I126").Value = 20490
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 61470
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -59000
$ws.Range("N126").Value = -13940

$ws.Range("H133").Value = 43835
$ws.Range("J133").Value = 43835
$ws.Range("L133").Value = 43835
$ws.Range("N133").Value = -53955

$ws.Range("H135").Value = 49664
$ws.Range("J135").Value = 49664
$ws.Range("L135").Value = 49664
$ws.Range("N135").Value = -59804

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2904.4614
$ws.Range("I7").Value = 2064.8
$ws.Range("J7").Value = 5703.3335
$ws.Range("K7").Value = 2064.8
$ws.Range("L7").Value = 5703.3335
$ws.Range("M7").Value = -1952.8
$ws.Range("N7").Value = -5927.3335

$ws.Range("H122").Value = 6270616.5
$ws.Range("I122").Value = 7151061.5
$ws.Range("J122").Value = 3335800
$ws.Range("K122").Value = 21453184.5
$ws.Range("L122").Value = 10007400
$ws.Range("M122").Value = -21450734.5
$ws.Range("N122").Value = -10012300

$ws.Range("H126").Value = 2904.4614
$ws.Range("I126").Value = 2064.8
$ws.Range("J126").Value = 5703.3335
$ws.Range("K126").Value = 6194.400000000001
$ws.Range("L126").Value = 17110.0005
$ws.Range("M126").Value = -3724.400000000001
$ws.Range("N126").Value = -22050.0005

$ws.Range("H133").Value = 85308.664
$ws.Range("J133").Value = 85308.664
$ws.Range("L133").Value = 85308.664
$ws.Range("N133").Value = -90368.664

$ws.Range("H138").Value = 47825.6
$ws.Range("J138").Value = 47825.6
$ws.Range("L138").Value = 47825.6
$ws.Range("N138").Value = -58105.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 37582.4
$ws.Range("J46").Value = 37582.4
$ws.Range("L46").Value = 37582.4
$ws.Range("N46").Value = -38044.4

$ws.Range("H132").Value = 2192.8928
$ws.Range("I132").Value = 966.05884
$ws.Range("J132").Value = 4088.9092
$ws.Range("K132").Value = 2898.17652
$ws.Range("L132").Value = 12266.7276
$ws.Range("M132").Value = -368.17652
$ws.Range("N132").Value = -17326.7276

$ws.Range("H134").Value = 37582.4
$ws.Range("J134").Value = 37582.4
$ws.Range("L134").Value = 112747.2
$ws.Range("N134").Value = -117817.2

$ws.Range("H136").Value = 2690979.2
$ws.Range("I136").Value = 2857.9656
$ws.Range("J136").Value = 5053267.5
$ws.Range("K136").Value = 8573.8968
$ws.Range("L136").Value = 15159802.5
$ws.Range("M136").Value = -6023.8968
$ws.Range("N136").Value = -15164902.5
